# Update market-price / profit figures on several sheets.
# Each sheet row holds currentAveragePrice(H)/NQ(I)/HQ(J), LevePrice NQ(K)/HQ(L)
# and the derived LeveProfit NQ(M)/HQ(N); this mirrors the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 11752.25
$ws.Range("I10").Value = 8499.5
$ws.Range("K10").Value = 8499.5
$ws.Range("M10").Value = -8206.5
$ws.Range("H112").Value = 1729.0769
$ws.Range("J112").Value = 1735.04
$ws.Range("L112").Value = 5205.12
$ws.Range("N112").Value = -7421.12
$ws.Range("H132").Value = 1133.6364
$ws.Range("I132").Value = 1133.6364
$ws.Range("K132").Value = 3400.9092
$ws.Range("M132").Value = -870.9092000000001
$ws.Range("H137").Value = 13701544
$ws.Range("I137").Value = 66669244
$ws.Range("K137").Value = 200007732
$ws.Range("M137").Value = -200005182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 18999.666
$ws.Range("J2").Value = 50500
$ws.Range("L2").Value = 50500
$ws.Range("N2").Value = -50726
$ws.Range("H32").Value = 5771.3706
$ws.Range("J32").Value = 26599.8
$ws.Range("L32").Value = 26599.8
$ws.Range("N32").Value = -27173.8
$ws.Range("H45").Value = 3353.1428
$ws.Range("I45").Value = 3245.3333
$ws.Range("K45").Value = 3245.3333
$ws.Range("M45").Value = -2868.3333
$ws.Range("H55").Value = 8500
$ws.Range("I55").Value = 8500
$ws.Range("K55").Value = 8500
$ws.Range("M55").Value = -8185
$ws.Range("H61").Value = 4882.88
$ws.Range("I61").Value = 4152.9
$ws.Range("K61").Value = 4152.9
$ws.Range("M61").Value = -3940.9
$ws.Range("H63").Value = 5411.9
$ws.Range("I63").Value = 3172.1538
$ws.Range("J63").Value = 9571.429
$ws.Range("K63").Value = 3172.1538
$ws.Range("L63").Value = 9571.429
$ws.Range("M63").Value = -2486.1538
$ws.Range("N63").Value = -10943.429
$ws.Range("H66").Value = 5411.9
$ws.Range("I66").Value = 3172.1538
$ws.Range("J66").Value = 9571.429
$ws.Range("K66").Value = 15860.769
$ws.Range("L66").Value = 47857.145
$ws.Range("M66").Value = -12428.769
$ws.Range("N66").Value = -54721.145
$ws.Range("H74").Value = 3181.7917
$ws.Range("I74").Value = 2139.4736
$ws.Range("J74").Value = 7142.6
$ws.Range("K74").Value = 2139.4736
$ws.Range("L74").Value = 7142.6
$ws.Range("M74").Value = -1265.4736
$ws.Range("N74").Value = -8890.6
$ws.Range("H77").Value = 3181.7917
$ws.Range("I77").Value = 2139.4736
$ws.Range("J77").Value = 7142.6
$ws.Range("K77").Value = 10697.368
$ws.Range("L77").Value = 35713
$ws.Range("M77").Value = -6329.367999999999
$ws.Range("N77").Value = -44449
$ws.Range("H80").Value = 133332
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 133332
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 133332
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -135328
$ws.Range("H83").Value = 133332
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 133332
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 399996
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -409980
$ws.Range("H97").Value = 1354.0667
$ws.Range("I97").Value = 1001.9091
$ws.Range("J97").Value = 2322.5
$ws.Range("K97").Value = 1001.9091
$ws.Range("L97").Value = 2322.5
$ws.Range("M97").Value = -505.9091
$ws.Range("N97").Value = -3314.5
$ws.Range("H102").Value = 3950
$ws.Range("I102").Value = 3950
$ws.Range("K102").Value = 3950
$ws.Range("M102").Value = -2328
$ws.Range("H116").Value = 18999.666
$ws.Range("J116").Value = 50500
$ws.Range("L116").Value = 50500
$ws.Range("N116").Value = -55088
$ws.Range("H132").Value = 2846.4783
$ws.Range("I132").Value = 1929.2858
$ws.Range("K132").Value = 5787.857400000001
$ws.Range("M132").Value = -3257.857400000001
$ws.Range("H136").Value = 4882.88
$ws.Range("I136").Value = 4152.9
$ws.Range("K136").Value = 12458.7
$ws.Range("M136").Value = -9908.699999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 18999.666
$ws.Range("J3").Value = 50500
$ws.Range("L3").Value = 50500
$ws.Range("N3").Value = -50728
$ws.Range("H54").Value = 1316.5
$ws.Range("I54").Value = 1316.5
$ws.Range("K54").Value = 1316.5
$ws.Range("M54").Value = -832.5
$ws.Range("H99").Value = 2326.25
$ws.Range("I99").Value = 2055.4
$ws.Range("K99").Value = 2055.4
$ws.Range("M99").Value = -557.4000000000001
$ws.Range("H107").Value = 706.1
$ws.Range("I107").Value = 706.1
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 706.1
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1213.9
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 2024.7407
$ws.Range("I134").Value = 1255.6342
$ws.Range("K134").Value = 3766.9026
$ws.Range("M134").Value = -1231.9026

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 83724120
$ws.Range("I4").Value = 143506860
$ws.Range("J4").Value = 15400996
$ws.Range("K4").Value = 430520580
$ws.Range("L4").Value = 46202988
$ws.Range("M4").Value = -430520468
$ws.Range("N4").Value = -46203212
$ws.Range("H107").Value = 18518804
$ws.Range("I107").Value = 222.81818
$ws.Range("J107").Value = 47619430
$ws.Range("K107").Value = 668.4545400000001
$ws.Range("L107").Value = 142858290
$ws.Range("M107").Value = 1251.54546
$ws.Range("N107").Value = -142862130
$ws.Range("H113").Value = 111113120
$ws.Range("J113").Value = 166668690
$ws.Range("L113").Value = 500006070
$ws.Range("N113").Value = -500010410
$ws.Range("H129").Value = 16674430
$ws.Range("J129").Value = 20842588
$ws.Range("L129").Value = 62527764
$ws.Range("N129").Value = -62537764
$ws.Range("H131").Value = 10306077
$ws.Range("J131").Value = 12154527
$ws.Range("L131").Value = 36463581
$ws.Range("N131").Value = -36473661
$ws.Range("H137").Value = 2196.16
$ws.Range("I137").Value = 983.8
$ws.Range("J137").Value = 4014.7
$ws.Range("K137").Value = 2951.4
$ws.Range("L137").Value = 12044.1
$ws.Range("M137").Value = 2148.6
$ws.Range("N137").Value = -22244.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 505533.5
$ws.Range("I80").Value = 1002386.2
$ws.Range("J80").Value = 8680.799999999999
$ws.Range("K80").Value = 1002386.2
$ws.Range("L80").Value = 8680.799999999999
$ws.Range("M80").Value = -1001388.2
$ws.Range("N80").Value = -10676.8
$ws.Range("H83").Value = 505533.5
$ws.Range("I83").Value = 1002386.2
$ws.Range("J83").Value = 8680.799999999999
$ws.Range("K83").Value = 5011931
$ws.Range("L83").Value = 43404
$ws.Range("M83").Value = -5006939
$ws.Range("N83").Value = -53388
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H126").Value = 4542.1763
$ws.Range("I126").Value = 2612.5715
$ws.Range("K126").Value = 7837.7145
$ws.Range("M126").Value = -5367.7145
$ws.Range("H132").Value = 2704.3096
$ws.Range("I132").Value = 2158.4167
$ws.Range("K132").Value = 6475.250100000001
$ws.Range("M132").Value = -3945.250100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9714.857
$ws.Range("I22").Value = 1750
$ws.Range("K22").Value = 1750
$ws.Range("M22").Value = -1455
$ws.Range("H27").Value = 9714.857
$ws.Range("I27").Value = 1750
$ws.Range("K27").Value = 1750
$ws.Range("M27").Value = -1643
$ws.Range("H39").Value = 32495
$ws.Range("J39").Value = 32495
$ws.Range("L39").Value = 32495
$ws.Range("N39").Value = -33415
$ws.Range("H55").Value = 3127870.5
$ws.Range("I55").Value = 5556050
$ws.Range("J55").Value = 5925.143
$ws.Range("K55").Value = 5556050
$ws.Range("L55").Value = 5925.143
$ws.Range("M55").Value = -5555877
$ws.Range("N55").Value = -6271.143
$ws.Range("H100").Value = 50011750
$ws.Range("I100").Value = 100005000
$ws.Range("K100").Value = 100005000
$ws.Range("M100").Value = -100004459
$ws.Range("H119").Value = 68421
$ws.Range("J119").Value = 68421
$ws.Range("L119").Value = 68421
$ws.Range("N119").Value = -78097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 9999.5
$ws.Range("I48").Value = 9999
$ws.Range("K48").Value = 9999
$ws.Range("M48").Value = -9430
$ws.Range("H75").Value = 49996.332
$ws.Range("I75").Value = 49996.332
$ws.Range("K75").Value = 49996.332
$ws.Range("M75").Value = -49060.332
$ws.Range("H78").Value = 49996.332
$ws.Range("I78").Value = 49996.332
$ws.Range("K78").Value = 149988.996
$ws.Range("M78").Value = -145308.996
$ws.Range("H119").Value = 75630.664
$ws.Range("J119").Value = 75630.664
$ws.Range("L119").Value = 75630.664
$ws.Range("N119").Value = -85306.664
$ws.Range("H126").Value = 3835
$ws.Range("I126").Value = 4365.8667
$ws.Range("K126").Value = 13097.6001
$ws.Range("M126").Value = -10607.6001
